# Generate Report for Handback
#
# The "505d041b-0c10-4e8c-afde-d41eb890b5b2.md" file has now been handed
# back (its Handback status flips from "Ready for handoff" /
# "not the latest" error to "Handed back: in sync with en-US", with fresh
# handback timestamps). As a result it now sorts to the first data row on
# every sheet, pushing the other two rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(2,1).Value = "505d041b-0c10-4e8c-afde-d41eb890b5b2.md"
$ov.Cells.Item(2,2).Value = "e2e\505d041b-0c10-4e8c-afde-d41eb890b5b2.md"
$ov.Cells.Item(2,3).Value = ".md"
$ov.Cells.Item(2,4).Value = ""
$ov.Cells.Item(2,5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(2,6).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(2,7).Value = "2016-07-26 08:20:40"

$ov.Cells.Item(3,1).Value = "fffff37654ce-3e64-4e5e-97fa-e71480b4877b.md"
$ov.Cells.Item(3,2).Value = "e2e\fffff37654ce-3e64-4e5e-97fa-e71480b4877b.md"
$ov.Cells.Item(3,3).Value = ".md"
$ov.Cells.Item(3,4).Value = ""
$ov.Cells.Item(3,5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(3,6).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(3,7).Value = "2016-07-26 08:09:56"

$ov.Cells.Item(4,1).Value = "ffffff73458b5e-c28f-4c01-9120-ffda3c258ae7.md"
$ov.Cells.Item(4,2).Value = "e2e\ffffff73458b5e-c28f-4c01-9120-ffda3c258ae7.md"
$ov.Cells.Item(4,3).Value = ".md"
$ov.Cells.Item(4,4).Value = ""
$ov.Cells.Item(4,5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(4,6).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(4,7).Value = "2016-07-26 08:09:56"

# rebuild hyperlinks (column B) to match new row order
while ($ov.Hyperlinks.Count -gt 0) {
    foreach ($h in $ov.Hyperlinks) { $h.Delete(); break }
}
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/b6985c375a33e082d1943053e3553409c92c6197/e2e/fffff37654ce-3e64-4e5e-97fa-e71480b4877b.md", "", "", "e2e\505d041b-0c10-4e8c-afde-d41eb890b5b2.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/18a71b264796044fae88c27d84ec4c93e5e2c444/e2e/ffffff73458b5e-c28f-4c01-9120-ffda3c258ae7.md", "", "", "e2e\fffff37654ce-3e64-4e5e-97fa-e71480b4877b.md")
$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/18a71b264796044fae88c27d84ec4c93e5e2c444/e2e/505d041b-0c10-4e8c-afde-d41eb890b5b2.md", "", "", "e2e\ffffff73458b5e-c28f-4c01-9120-ffda3c258ae7.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Cells.Item(2,1).Value  = "505d041b-0c10-4e8c-afde-d41eb890b5b2.md"
$zh.Cells.Item(2,2).Value  = ".md"
$zh.Cells.Item(2,3).Value  = "Handed back: in sync with en-US"
$zh.Cells.Item(2,4).Value  = "e2e"
$zh.Cells.Item(2,5).Value  = "ht"
$zh.Cells.Item(2,6).Value  = "505d041b-0c10-4e8c-afde-d41eb890b5b2.de55513c87b37a3b409a404a7b52f4aec8d61f69.zh-cn.xlf"
$zh.Cells.Item(2,7).Value  = "2016-07-26 08:20:30"
$zh.Cells.Item(2,8).Value  = "505d041b-0c10-4e8c-afde-d41eb890b5b2.md"
$zh.Cells.Item(2,9).Value  = "505d041b-0c10-4e8c-afde-d41eb890b5b2.de55513c87b37a3b409a404a7b52f4aec8d61f69.zh-cn.xlf"
$zh.Cells.Item(2,10).Value = "2016-07-26 08:21:15"
$zh.Cells.Item(2,11).Value = ""
$zh.Cells.Item(2,12).Value = "True"
$zh.Cells.Item(2,13).Value = ""
$zh.Cells.Item(2,14).Value = "False"
$zh.Cells.Item(2,15).Value = ""

$zh.Cells.Item(3,1).Value  = "fffff37654ce-3e64-4e5e-97fa-e71480b4877b.md"
$zh.Cells.Item(3,2).Value  = ".md"
$zh.Cells.Item(3,3).Value  = "Handed back: in sync with en-US"
$zh.Cells.Item(3,4).Value  = "e2e"
$zh.Cells.Item(3,5).Value  = "ht"
$zh.Cells.Item(3,6).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.zh-cn.xlf"
$zh.Cells.Item(3,7).Value  = "2016-07-26 08:09:46"
$zh.Cells.Item(3,8).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md"
$zh.Cells.Item(3,9).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.zh-cn.xlf"
$zh.Cells.Item(3,10).Value = "2016-07-26 08:10:36"
$zh.Cells.Item(3,11).Value = ""
$zh.Cells.Item(3,12).Value = "True"
$zh.Cells.Item(3,13).Value = ""
$zh.Cells.Item(3,14).Value = "False"
$zh.Cells.Item(3,15).Value = ""

$zh.Cells.Item(4,1).Value  = "ffffff73458b5e-c28f-4c01-9120-ffda3c258ae7.md"
$zh.Cells.Item(4,2).Value  = ".md"
$zh.Cells.Item(4,3).Value  = "Handed back: in sync with en-US"
$zh.Cells.Item(4,4).Value  = "e2e"
$zh.Cells.Item(4,5).Value  = "ht"
$zh.Cells.Item(4,6).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.zh-cn.xlf"
$zh.Cells.Item(4,7).Value  = "2016-07-26 08:09:46"
$zh.Cells.Item(4,8).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md"
$zh.Cells.Item(4,9).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.zh-cn.xlf"
$zh.Cells.Item(4,10).Value = "2016-07-26 08:10:36"
$zh.Cells.Item(4,11).Value = ""
$zh.Cells.Item(4,12).Value = "True"
$zh.Cells.Item(4,13).Value = ""
$zh.Cells.Item(4,14).Value = "False"
$zh.Cells.Item(4,15).Value = ""

while ($zh.Hyperlinks.Count -gt 0) {
    foreach ($h in $zh.Hyperlinks) { $h.Delete(); break }
}
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/b6985c375a33e082d1943053e3553409c92c6197/e2e/fffff37654ce-3e64-4e5e-97fa-e71480b4877b.md", "", "", "505d041b-0c10-4e8c-afde-d41eb890b5b2.md")
$zh.Hyperlinks.Add($zh.Range("H2"), "https://github.com/OpenLocalizationTestOrg/ol-test-zhcn/blob/5ad3daaf78e09463c29f71bdb7603bd48d532cc0/e2e/dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md", "", "", "505d041b-0c10-4e8c-afde-d41eb890b5b2.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/18a71b264796044fae88c27d84ec4c93e5e2c444/e2e/ffffff73458b5e-c28f-4c01-9120-ffda3c258ae7.md", "", "", "fffff37654ce-3e64-4e5e-97fa-e71480b4877b.md")
$zh.Hyperlinks.Add($zh.Range("H3"), "https://github.com/OpenLocalizationTestOrg/ol-test-zhcn/blob/5ad3daaf78e09463c29f71bdb7603bd48d532cc0/e2e/dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md", "", "", "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/18a71b264796044fae88c27d84ec4c93e5e2c444/e2e/505d041b-0c10-4e8c-afde-d41eb890b5b2.md", "", "", "ffffff73458b5e-c28f-4c01-9120-ffda3c258ae7.md")
$zh.Hyperlinks.Add($zh.Range("H4"), "https://github.com/OpenLocalizationTestOrg/ol-test-zhcn/blob/7f04942f2976ceb293805dab206e79f5c4e0c7d6/e2e/505d041b-0c10-4e8c-afde-d41eb890b5b2.md", "", "", "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md")

# Error Detail column (O) is now empty on every row -> narrower column
$zh.Columns.Item(15).ColumnWidth = 12.8

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Cells.Item(2,1).Value  = "505d041b-0c10-4e8c-afde-d41eb890b5b2.md"
$de.Cells.Item(2,2).Value  = ".md"
$de.Cells.Item(2,3).Value  = "Handed back: in sync with en-US"
$de.Cells.Item(2,4).Value  = "e2e"
$de.Cells.Item(2,5).Value  = "ht"
$de.Cells.Item(2,6).Value  = "505d041b-0c10-4e8c-afde-d41eb890b5b2.de55513c87b37a3b409a404a7b52f4aec8d61f69.de-de.xlf"
$de.Cells.Item(2,7).Value  = "2016-07-26 08:20:40"
$de.Cells.Item(2,8).Value  = "505d041b-0c10-4e8c-afde-d41eb890b5b2.md"
$de.Cells.Item(2,9).Value  = "505d041b-0c10-4e8c-afde-d41eb890b5b2.de55513c87b37a3b409a404a7b52f4aec8d61f69.de-de.xlf"
$de.Cells.Item(2,10).Value = "2016-07-26 08:21:29"
$de.Cells.Item(2,11).Value = ""
$de.Cells.Item(2,12).Value = "True"
$de.Cells.Item(2,13).Value = ""
$de.Cells.Item(2,14).Value = "False"
$de.Cells.Item(2,15).Value = ""

$de.Cells.Item(3,1).Value  = "fffff37654ce-3e64-4e5e-97fa-e71480b4877b.md"
$de.Cells.Item(3,2).Value  = ".md"
$de.Cells.Item(3,3).Value  = "Handed back: in sync with en-US"
$de.Cells.Item(3,4).Value  = "e2e"
$de.Cells.Item(3,5).Value  = "ht"
$de.Cells.Item(3,6).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.de-de.xlf"
$de.Cells.Item(3,7).Value  = "2016-07-26 08:09:56"
$de.Cells.Item(3,8).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md"
$de.Cells.Item(3,9).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.de-de.xlf"
$de.Cells.Item(3,10).Value = "2016-07-26 08:10:51"
$de.Cells.Item(3,11).Value = ""
$de.Cells.Item(3,12).Value = "True"
$de.Cells.Item(3,13).Value = ""
$de.Cells.Item(3,14).Value = "False"
$de.Cells.Item(3,15).Value = ""

$de.Cells.Item(4,1).Value  = "ffffff73458b5e-c28f-4c01-9120-ffda3c258ae7.md"
$de.Cells.Item(4,2).Value  = ".md"
$de.Cells.Item(4,3).Value  = "Handed back: in sync with en-US"
$de.Cells.Item(4,4).Value  = "e2e"
$de.Cells.Item(4,5).Value  = "ht"
$de.Cells.Item(4,6).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.de-de.xlf"
$de.Cells.Item(4,7).Value  = "2016-07-26 08:09:56"
$de.Cells.Item(4,8).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md"
$de.Cells.Item(4,9).Value  = "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.2f07c35368b9579b291927cc6804ae3b6a7af3f0.de-de.xlf"
$de.Cells.Item(4,10).Value = "2016-07-26 08:10:51"
$de.Cells.Item(4,11).Value = ""
$de.Cells.Item(4,12).Value = "True"
$de.Cells.Item(4,13).Value = ""
$de.Cells.Item(4,14).Value = "False"
$de.Cells.Item(4,15).Value = ""

while ($de.Hyperlinks.Count -gt 0) {
    foreach ($h in $de.Hyperlinks) { $h.Delete(); break }
}
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/b6985c375a33e082d1943053e3553409c92c6197/e2e/fffff37654ce-3e64-4e5e-97fa-e71480b4877b.md", "", "", "505d041b-0c10-4e8c-afde-d41eb890b5b2.md")
$de.Hyperlinks.Add($de.Range("H2"), "https://github.com/OpenLocalizationTestOrg/ol-test-dede/blob/f6563957604f1d0d741a5f9b715f90017cb1fc2e/e2e/dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md", "", "", "505d041b-0c10-4e8c-afde-d41eb890b5b2.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/18a71b264796044fae88c27d84ec4c93e5e2c444/e2e/ffffff73458b5e-c28f-4c01-9120-ffda3c258ae7.md", "", "", "fffff37654ce-3e64-4e5e-97fa-e71480b4877b.md")
$de.Hyperlinks.Add($de.Range("H3"), "https://github.com/OpenLocalizationTestOrg/ol-test-dede/blob/f6563957604f1d0d741a5f9b715f90017cb1fc2e/e2e/dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md", "", "", "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/18a71b264796044fae88c27d84ec4c93e5e2c444/e2e/505d041b-0c10-4e8c-afde-d41eb890b5b2.md", "", "", "ffffff73458b5e-c28f-4c01-9120-ffda3c258ae7.md")
$de.Hyperlinks.Add($de.Range("H4"), "https://github.com/OpenLocalizationTestOrg/ol-test-dede/blob/a78ad420daaa2d83916831a77b764a19dd34b6b3/e2e/505d041b-0c10-4e8c-afde-d41eb890b5b2.md", "", "", "dfd87fc3-e78e-4fa1-95a1-f875ed042ef9.md")

# Error Detail column (O) is now empty on every row -> narrower column
$de.Columns.Item(15).ColumnWidth = 12.8
